# Updates based on progress report 2 feedback
# Rewrites the two summary tables (RQ1 cause-of-flakiness / RQ2 fix-for-flakiness
# and programming-language breakdown) with the revised taxonomy & counts, resizes
# the three tables accordingly, turns on the totals row for the "fix" table, and
# moves the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clear the old data area (rows that shrink / disappear) so nothing is left
#    behind once the ranges are resized.
# ---------------------------------------------------------------------------
$ws.Range("B3:F30").ClearContents()

# ---------------------------------------------------------------------------
# 2. Write the new "RQ1: Cause of Flakiness?" / "RQ2: Fix for Flakiness?" table
#    headers (row 3, unchanged values but rewritten for completeness).
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "RQ1: Cause of Flakiness?"
$ws.Range("C3").Value = "Count"
$ws.Range("E3").Value = "RQ2: Fix for Flakiness?"
$ws.Range("F3").Value = "Count"

# ---------------------------------------------------------------------------
# 3. RQ1: Cause of Flakiness (Table1, B4:C14)
# ---------------------------------------------------------------------------
$causeRows = @(
    @("Memory", 6),
    @("Dependencies / Environment", 10),
    @("Order of Events", 23),
    @("Concurrency", 6),
    @("Async Wait", 10),
    @("Delay", 7),
    @("Collections", 3),
    @("Bit Manipulation / Arithmetic", 3),
    @("Comparisons", 1),
    @("Algorithmic Flakiness", 14),
    @("Incorrect Logic", 13)
)
$r = 4
foreach ($row in $causeRows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. RQ2: Fix for Flakiness (Table13, E4:F22, with totals row on row 23)
# ---------------------------------------------------------------------------
$fixRows = @(
    @("relax acceptance/tolerance", 15),
    @("deallocate unused space", 5),
    @("update ", 3),
    @("implement event ordering", 8),
    @("refactor program logic", 13),
    @("setup state", 4),
    @("tear down state after shutdown", 11),
    @("locks", 5),
    @("added waitFor", 9),
    @("added memory for test", 1),
    @("add custom delay / wait", 7),
    @("sort / establish ordering", 2),
    @("limit elements", 1),
    @("bit clear", 1),
    @("update global state", 1),
    @("Promise statement", 1),
    @("remove dependencies", 6),
    @("same type", 1),
    @("no hardcoded values", 2)
)
$r = 4
foreach ($row in $fixRows) {
    $ws.Cells.Item($r, 5).Value = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5. Programming Language table (Table3), header now on row 17, data 18:29
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = "Programming Language"
$ws.Range("C17").Value = "Count"

# NOTE: Table3's header row is moving up one row (18 -> 17). When ListObject.Resize
# shifts the header row, the engine re-derives the tableColumn names from whatever
# currently sits at the *old* header location (row 18) rather than the new one.
# Stash the intended header text there too so the resize captures the right
# names; the real data for row 18 ("Scala", 3) is written right afterwards.
$ws.Range("B18").Value = "Programming Language"
$ws.Range("C18").Value = "Count"

$loLang = $ws.ListObjects.Item(3)
$loLang.Resize($ws.Range("B17:C29"))

$langRows = @(
    @("Scala", 3),
    @("Swift", 10),
    @("TypeScript", 8),
    @("Python", 10),
    @("JavaScript", 5),
    @("Java", 9),
    @("C++", 19),
    @("C#", 10),
    @("Go", 9),
    @("Kotlin", 8),
    @("Clojure", 4),
    @("Rust", 1)
)
$r = 18
foreach ($row in $langRows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 6. Resize the remaining tables (ListObjects) to reflect the new data extents.
#    (Their header rows do not move, so no workaround is needed for them.)
# ---------------------------------------------------------------------------
$loCause = $ws.ListObjects.Item(1)
$loCause.Resize($ws.Range("B3:C15"))

$loFix = $ws.ListObjects.Item(2)
$loFix.Resize($ws.Range("E3:F22"))
$loFix.ShowTotals = $true

# ---------------------------------------------------------------------------
# 7. Update the selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("F23").Select()
